$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix ADDX bug: O6 (Ain) should be 1, not 0 ---
$ws.Range("O6").Value = "1"

# --- Add new SUBX row (row 7) ---
# Copy formatting (number format/style) from row 6's instruction columns (A:P)
# down to row 7 so the new row matches the sheet's existing look (text-formatted,
# centered cells) without touching column Q, which already holds the shared
# CONCATENATE formula for every row down to 53.
$ws.Range("A6:P6").Copy()
$ws.Range("A7:P7").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Fill in the SUBX instruction data
$ws.Range("A7").Value = "SUBX;a-x->a"
$ws.Range("B7").Value = "01"
$ws.Range("C7").Value = "100"
$ws.Range("D7").Value = "1"
$ws.Range("E7").Value = "0"
$ws.Range("F7").Value = "1"
$ws.Range("G7").Value = "0"
$ws.Range("H7").Value = "0"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "1"
$ws.Range("L7").Value = "0"
$ws.Range("M7").Value = "1"
$ws.Range("N7").Value = "0"
$ws.Range("O7").Value = "1"
$ws.Range("P7").Value = "0010010"

# Move the active selection to Q7, matching the author's final cursor position
$ws.Range("Q7").Select() | Out-Null
